$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Redesign: add a new "Primer Apellido" (first surname) column in C.
# The existing "Segundo Apellido" (second surname) column stays in D,
# and columns D:H keep their positions/values unchanged.
$ws.Range("C1").Value = "Primer Apellido"
$ws.Range("C2").Value = "Pérez"
$ws.Range("C3").Value = "Ramírez"
$ws.Range("C4").Value = "de la rosa"
$ws.Range("C5").Value = "Magdalena"
$ws.Range("C6").Value = "Annia"
